# In the "Output" section, two paragraphs are removed:
#   - the bordered paragraph containing "Ingen." (with the _GoBack bookmark)
#   - the empty paragraph that followed it
# so that "Output" is immediately followed by "Ingen dokumentlyter er registrert."
$d = $word.ActiveDocument

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Ingen.") {
        $next = $p.Next()
        if ($next -ne $null -and $next.Range.Text.TrimEnd([char]13, [char]7) -eq "") {
            $after = $next.Next()
            if ($after -ne $null -and $after.Range.Text.StartsWith("Ingen dokumentlyter")) {
                $target = $p
            }
        }
    }
}

if ($target -ne $null) {
    $startPara = $target
    $endPara = $target.Next()
    $rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $rng.Delete()
}
